$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Cell content edits -----------------------------------------------
# Rows where the "heriot assessed" (J) answer changed from "No" to
# "Unknown", and the now-meaningless "heriot amt" (K) value of 0 was
# cleared out.
$rowsNoToUnknown = @(3, 7, 48, 52, 57)
foreach ($r in $rowsNoToUnknown) {
    $ws.Range("J$r").Value = "Unknown"
    $ws.Range("K$r").ClearContents()
}

# "--" placeholders in the heriot animal column (L) removed.
$ws.Range("L16").ClearContents()
$ws.Range("L17").ClearContents()

# Stray "N/A" heriot amount cleared.
$ws.Range("K63").ClearContents()

# --- Column width -------------------------------------------------------
# Column L (heriot animal) gets an explicit width.
$ws.Columns.Item(12).ColumnWidth = 5

# --- View / selection state ---------------------------------------------
# Scroll back to the top of the frozen pane and select A7.
[void]$ws.Range("A7").Select()
